# Updated cryptos list on Sun Oct  6 20:40:10 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) columns of the
# cryptocurrency table with the latest scraped values. Numeric-looking
# price strings are forced to remain text (matching the source feed's
# formatting, e.g. trailing zeros) by briefly applying a text number
# format before the write and clearing the format immediately after so
# no residual cell styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.635.76"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "2.437.66"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.42"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.31"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.86"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.65%  "
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "62.420.72"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").Value = "2.435.80"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.97"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.85"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.35"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "580.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("E27").Value = "  +9.33%  "
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.42"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.71%  "
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.76"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "148.11"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("E41").Value = "  +2.58%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +10.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.25"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.56"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.23%  "
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0231"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  +5.17%  "
